# Updated cryptos list on Sat Jul 22 20:00:06 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.
# D-column values are plain text that merely look numeric (e.g. "1.892.72"),
# so they are written via a leading apostrophe (Formula = "'<text>") to force
# text storage, then the style is reset to "Normal" so no residual
# quote-prefix/number-format is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'29.845.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Formula = "'1.887.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Formula = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Formula = "'0.7485"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.59%  "
$ws.Range("D6").Formula = "'242.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Formula = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Formula = "'0.3119"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Formula = "'25.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Formula = "'0.07125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").Formula = "'0.08489"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.76%  "
$ws.Range("D12").Formula = "'0.7607"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Formula = "'1.892.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Formula = "'5.357"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Formula = "'93.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Formula = "'6.145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Formula = "'29.867.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Formula = "'13.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Formula = "'243.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Formula = "'0.000007795"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Formula = "'2.156.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Formula = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Formula = "'7.986"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Formula = "'0.9998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Formula = "'9.353"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Formula = "'162.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Formula = "'18.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Formula = "'2.030"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Formula = "'1.495"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.75%  "
$ws.Range("D31").Formula = "'1.533"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Formula = "'4.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Formula = "'4.103"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Formula = "'0.05392"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").Formula = "'1.239"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Formula = "'0.7451"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Formula = "'1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Formula = "'2.704"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Formula = "'0.01934"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Formula = "'2.770"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Formula = "'0.4458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Formula = "'6.061"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Formula = "'1.090.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.69%  "
$ws.Range("D44").Formula = "'72.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Formula = "'0.8537"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Formula = "'102.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Formula = "'7.686"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").Formula = "'1.863"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").Formula = "'3.064"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Formula = "'2.052.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
